$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "\n" placeholder values stored in ExEn/ExVi for the "Bye" row were
# bogus leftovers — clear them out entirely (leader board topic-xp-stats
# cleanup), so the cells disappear rather than staying as empty strings.
$ws.Range("D3:E3").ClearContents()

# Move the saved selection from G7 to F7, matching the workbook's last
# recorded view state.
$ws.Range("F7").Select()
